# UPDATE: Inventory quantity adjustments and history tracking
#
# Adds two new inventory rows (id 41 "Y27632" and id 42 "E8 Supplement")
# to the bottom of the tracked table on Sheet1, then leaves the
# selection on H49 (matching the author's final cursor position).
#
# Cell-write order below intentionally mirrors the order new shared
# strings first appear in the saved workbook (name/unit columns first,
# then the datasheet URLs last) so new shared-string table entries come
# out in the same sequence as the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42 (id 41): Y27632 -------------------------------------------------
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = "Y27632"
$ws.Range("C42").Value = 72304
$ws.Range("E42").Value = "Small Molecule"
$ws.Range("F42").Value = "Aliquot"
$ws.Range("G42").Value = -30
$ws.Range("H42").Value = 9
$ws.Range("I42").Value = "Opened"
$ws.Range("J42").Value = 60
$ws.Range("K42").Value = "70 uL"

# --- Row 43 (id 42): E8 Supplement ------------------------------------------
$ws.Range("A43").Value = 42
$ws.Range("B43").Value = "E8 Supplement"
$ws.Range("C43").Value = "A1517-01"
$ws.Range("E43").Value = "Supplement"
$ws.Range("F43").Value = "Aliquot"
$ws.Range("G43").Value = -30
$ws.Range("H43").Value = 9
$ws.Range("I43").Value = "Opened"
$ws.Range("J43").Value = 11
$ws.Range("K43").Value = "1 mL"

# --- datasheet_url column filled in last for both new rows ------------------
$ws.Range("D42").Value = "https://www.stemcell.com/products/y-27632.html"
$ws.Range("D43").Value = "https://www.thermofisher.com/order/catalog/product/A1517001"

# --- restore the cursor/selection position left behind by the edit ---------
$ws.Range("H49").Select()
